# Update the two test e-mail addresses used by the Selenium/Chrome-headless
# test fixture (userData.xlsx) and mirror the workbook window-size change
# that Excel recorded when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C1: es234@yahoo.com -> es2345@yahoo.com
$ws.Range("C1").Value = "es2345@yahoo.com"

# C2: us1234@test.com -> us12345@test.com
$ws.Range("C2").Value = "us12345@test.com"

# Best-effort: record the new workbook window size (xWindow/yWindow stay 0).
# Harmless if the host does not persist window chrome metrics.
$win = $excel.ActiveWindow
$win.Width = 11430
$win.Height = 2670
